$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = 'Datos actualizados a 26 de Mayo de 2020 a las 00:35'
$ws.Range("B4").Value = 1704061
$ws.Range("C4").Value = 17625
$ws.Range("D4").Value = 461154
$ws.Range("E4").Value = 1143153
$ws.Range("B16").Value = 85698
$ws.Range("C16").Value = 999
$ws.Range("D16").Value = 44593
$ws.Range("E16").Value = 34564
$ws.Range("G16").Value = 117
$ws.Range("H16").Value = 6541
$ws.Range("A88").Value = 'Gabon'
$ws.Range("B88").Value = 2135
$ws.Range("C88").Value = 201
$ws.Range("D88").Value = 562
$ws.Range("E88").Value = 1559
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 14
$ws.Range("A89").Value = 'Republica de Macedonia'
$ws.Range("B89").Value = 1999
$ws.Range("C89").Value = 21
$ws.Range("D89").Value = 1439
$ws.Range("E89").Value = 447
$ws.Range("H89").Value = 113
$ws.Range("A90").Value = 'El Salvador'
$ws.Range("B90").Value = 1983
$ws.Range("C90").Value = 68
$ws.Range("D90").Value = 698
$ws.Range("E90").Value = 1250
$ws.Range("H90").Value = 35
$ws.Range("A91").Value = 'Cuba'
$ws.Range("B91").Value = 1947
$ws.Range("C91").Value = 6
$ws.Range("D91").Value = 1704
$ws.Range("E91").Value = 161
$ws.Range("H91").Value = 82
$ws.Range("B164").Value = 137
$ws.Range("C164").Value = 2
$ws.Range("E164").Value = 64
$ws.Range("G164").Value = 1
$ws.Range("H164").Value = 11
$ws.Range("A197").Value = 'Fiyi'
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 15
$ws.Range("H197").Value = 0
$ws.Range("A198").Value = 'Curazao'
$ws.Range("C198").Value = 1
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1
$ws.Range("A199").Value = 'Santa Lucia'
$ws.Range("A201").Value = 'Nueva Caledonia'
$ws.Range("A215").Value = 'San Bartolome'
$ws.Range("A216").Value = 'Bonaire, San Eustaquio y Saba'

Write-Output "Applied 53 cell updates"
